# Apply the additive edits described by the diff: fill in previously-empty
# "carrier" (D) and "pair_kind" (J) columns for the practice/generic rows,
# and populate the new unique_video / unique_audio block (rows 14-21,
# columns C/D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (2-5): carrier column D gets the matching "can/where/do/look"
# word that is already used as the pair_kind's corresponding generic carrier.
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic rows (6-9): new pair_kind column J marks the unique video/audio
# stimulus type for this pair.
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# New rows 14-21: kind (C) + carrier (D) for the unique_video / unique_audio
# stimulus entries, following the same carrier pattern (can, can, do, do,
# look, look, where, where) as the original generic rows.
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"
$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"
$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"
$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"
$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"
$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"
$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
